# Agregar ejemplo de producto inactivo (Filtro Prefiltro G4) en la fila 3
# de la hoja 'Productos', replicando el layout de la fila 2 (producto activo).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Las columnas con valores puramente numéricos deben forzarse a texto
# (formato '@') para que Excel no las reinterprete como números y así
# conserven el mismo formato que las demás celdas de datos (texto).
$ws.Range('E3').NumberFormat = '@'
$ws.Range('G3').NumberFormat = '@'
$ws.Range('K3').NumberFormat = '@'

$ws.Range('A3').Value = 'Filtro Prefiltro G4'
$ws.Range('B3').Value = 'Pre-filter G4'
$ws.Range('C3').Value = 'Prefiltro de baja eficiencia para sistemas de ventilación'
$ws.Range('D3').Value = 'Low efficiency pre-filter for ventilation systems'
$ws.Range('E3').Value = '250.00'
$ws.Range('F3').Value = 'MXN'
$ws.Range('G3').Value = '12.50'
$ws.Range('H3').Value = 'Filtros de Aire'
$ws.Range('I3').Value = 'Air Filters'
$ws.Range('J3').Value = 'Prefiltros'
$ws.Range('K3').Value = '50'
$ws.Range('L3').Value = 'inactive'
$ws.Range('M3').Value = 'prefiltro, ventilacion'
$ws.Range('N3').Value = '20x20x2'
$ws.Range('O3').Value = '20x20x2'
$ws.Range('P3').Value = '500x500x50mm'
$ws.Range('Q3').Value = '1.5kg'
$ws.Range('R3').Value = 'Cartón'
$ws.Range('S3').Value = '6 meses'
$ws.Range('U3').Value = 'PRE-G4-20-20'
$ws.Range('V3').Value = '35% a 0.3µm'
$ws.Range('W3').Value = '35% at 0.3µm'
$ws.Range('X3').Value = 'G4'
$ws.Range('Y3').Value = 'Marco de cartón, fácil instalación'
$ws.Range('Z3').Value = 'Cardboard frame, easy installation'
$ws.Range('AA3').Value = 'Cartón'
$ws.Range('AB3').Value = '60°C'
$ws.Range('AC3').Value = 'Sistemas de ventilación general'
$ws.Range('AD3').Value = 'General ventilation systems'
$ws.Range('AE3').Value = 'Oficinas, comercios'
$ws.Range('AF3').Value = 'Offices, commercial'
$ws.Range('AG3').Value = 'Bajo costo, fácil mantenimiento'
$ws.Range('AH3').Value = 'Low cost, easy maintenance'
